$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.049036860466003
$ws.Range("B1").Value = 2.809235334396362
$ws.Range("C1").Value = 8.820466995239258
$ws.Range("D1").Value = 2.035429239273071
$ws.Range("E1").Value = 1.139391183853149
